$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.969.48'
$ws.Range("E2").Value = '  +3.22%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.419.67'
$ws.Range("E3").Value = '  +2.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.37'
$ws.Range("E5").Value = '  +1.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.30'
$ws.Range("E6").Value = '  +2.33%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +3.38%  '

# Row 9
$ws.Range("E9").Value = '  +1.03%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("E10").Value = '  +4.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.357'
$ws.Range("E11").Value = '  +0.37%  '

# Row 12
$ws.Range("E12").Value = '  -1.81%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.83'
$ws.Range("E13").Value = '  +4.49%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.852.76'
$ws.Range("E14").Value = '  +2.92%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.918.03'
$ws.Range("E15").Value = '  +3.23%  '

# Row 16
$ws.Range("E16").Value = '  +1.37%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.398.55'
$ws.Range("E17").Value = '  +2.05%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.34'
$ws.Range("E18").Value = '  +4.97%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.38'
$ws.Range("E19").Value = '  +1.70%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '331.39'
$ws.Range("E20").Value = '  +0.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.73'
$ws.Range("E21").Value = '  -0.16%  '

# Row 22
$ws.Range("E22").Value = '  +0.09%  '

# Row 23
$ws.Range("E23").Value = '  +3.68%  '

# Row 24
$ws.Range("E24").Value = '  +3.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.62'
$ws.Range("E25").Value = '  +3.92%  '

# Row 26
$ws.Range("E26").Value = '  +0.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.34'
$ws.Range("E27").Value = '  +1.61%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0782'
$ws.Range("E28").Value = '  +6.18%  '

# Row 29
$ws.Range("E29").Value = '  +0.87%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.17'
$ws.Range("E30").Value = '  -0.08%  '

# Row 32
$ws.Range("E32").Value = '  +1.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.02'
$ws.Range("E33").Value = '  +1.48%  '

# Row 35
$ws.Range("E35").Value = '  +4.80%  '

# Row 36
$ws.Range("E36").Value = '  +0.15%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.18'
$ws.Range("E37").Value = '  -0.17%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("E38").Value = '  +1.09%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.62'
$ws.Range("E39").Value = '  +1.45%  '

# Row 40
$ws.Range("E40").Value = '  +9.37%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '313.02'
$ws.Range("E41").Value = '  +7.74%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.68'
$ws.Range("E42").Value = '  +0.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.31'
$ws.Range("E43").Value = '  -0.46%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0967'
$ws.Range("E44").Value = '  +1.75%  '

# Row 45
$ws.Range("E45").Value = '  +1.58%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.411'
$ws.Range("E46").Value = '  +7.72%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.576'
$ws.Range("E47").Value = '  +1.64%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.27'
$ws.Range("E48").Value = '  +1.74%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0225'
$ws.Range("E49").Value = '  +1.15%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.63'
$ws.Range("E50").Value = '  +0.76%  '

# Row 51
$ws.Range("E51").Value = '  -0.22%  '
